# "Generate Report for Handoff"
# Updates the localization-status report to reflect that the report is now
# ready for handoff: status cells flip from "In Translation" to
# "Ready for handoff", and the associated handoff timestamps are refreshed.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet -------------------------------------------------------
# E2 = zh-cn status, F2 = de-de status, G2 = latest handoff xliff generate date
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-18 23:00:07"

# --- zh-cn sheet ------------------------------------------------------------
# C2 = Status, H2 = Latest Handoff Datetime
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-18 22:59:56"

# --- de-de sheet ------------------------------------------------------------
# C2 = Status (H2's handoff datetime stays the same text as before, which now
# also matches the refreshed Overview "Latest HO Xliff Generate Date" value)
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-18 23:00:07"

# --- Column width refresh (status columns grew wider to fit new text) -------
$wsOverview.Columns.Item(5).ColumnWidth = 16.3333333333333
$wsOverview.Columns.Item(6).ColumnWidth = 16.3333333333333
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3333333333333
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3333333333333
